# FIX SYSTEM CONTENT IN api.py
# Update the FINE-TUNING sheet with new checkpoint/model data (llmv4 & llmv5),
# relabel the checkpoint sub-headers, add a token-price reference cell (Q2)
# and make the CLP cost formulas reference it instead of a hard-coded 996.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FINE-TUNING")

# --- New reference cell: price per 1K tokens (CLP), used by column Q formulas ---
$ws.Range("Q2").Value = 942

# --- Row 4: rename checkpoint sub-headers from step numbers to qualitative labels ---
$ws.Range("H4").Value = "LOW"
$ws.Range("I4").Value = "MEDIUM"
$ws.Range("J4").Value = "HIGH"

# --- Column Q (CLP cost): every row's formula now references $Q$2 instead of literal 996 ---
$ws.Range("Q5").Formula = "=P5*`$Q`$2"
$ws.Range("Q6").Formula = "=P6*`$Q`$2"
$ws.Range("Q7").Formula = "=P7*`$Q`$2"
$ws.Range("Q8").Formula = "=P8*`$Q`$2"
$ws.Range("Q9").Formula = "=P9*`$Q`$2"
$ws.Range("Q10").Formula = "=P10*`$Q`$2"

# --- Row 8 (model 4, dots-llmv4): fill in previously-empty checkpoint/training data ---
$ws.Range("H8").Value = "ft:gpt-4o-mini-2024-07-18:personal:dots-llmv4:B346Az6s:ckpt-step-8"
$ws.Range("I8").Value = "ft:gpt-4o-mini-2024-07-18:personal:dots-llmv4:B346B2Fp:ckpt-step-12"
$ws.Range("J8").Value = "ft:gpt-4o-mini-2024-07-18:personal:dots-llmv4:B346BIAP"
$ws.Range("K8").Value = 143632
$ws.Range("L8").Value = 4
$ws.Range("M8").Value = 0.1
$ws.Range("N8").Value = 4
$ws.Range("O8").Value = 1937114649
$ws.Rows.Item(8).RowHeight = 30

# --- Row 9 (model 5, dots-llmv5): fill in previously-empty checkpoint/training data ---
$ws.Range("H9").Value = "ft:gpt-4o-mini-2024-07-18:personal:dots-llmv5:B34rzjN4:ckpt-step-70"
$ws.Range("I9").Value = "ft:gpt-4o-mini-2024-07-18:personal:dots-llmv5:B34rzxJo:ckpt-step-84"
$ws.Range("J9").Value = "ft:gpt-4o-mini-2024-07-18:personal:dots-llmv5:B34s0MaY"
$ws.Range("K9").Value = 251356
$ws.Range("L9").Value = 7
$ws.Range("M9").Value = 1
$ws.Range("N9").Value = 1.8
$ws.Range("O9").Value = 2044383751
$ws.Rows.Item(9).RowHeight = 30

# --- Row 10 (model 6): clear the placeholder training params, no checkpoint data yet ---
$ws.Range("L10").ClearContents()
$ws.Range("M10").ClearContents()
$ws.Range("N10").ClearContents()
$ws.Rows.Item(10).RowHeight = 30

# --- Row 12: a stray formatted (underlined) blank cell left below the table ---
$ws.Range("L12").Font.Underline = 2

$wb.Save()
